$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1750.25
$ws.Range("I43").Value = 1499.5
$ws.Range("J43").Value = 2001
$ws.Range("K43").Value = 1499.5
$ws.Range("L43").Value = 2001
$ws.Range("M43").Value = -1430.5
$ws.Range("N43").Value = -2139
$ws.Range("H98").Value = 3240.5186
$ws.Range("I98").Value = 3194.3809
$ws.Range("J98").Value = 3402
$ws.Range("K98").Value = 3194.3809
$ws.Range("L98").Value = 3402
$ws.Range("M98").Value = -1696.3809
$ws.Range("N98").Value = -6398
$ws.Range("H116").Value = 5953.3335
$ws.Range("I116").Value = 3740.3333
$ws.Range("K116").Value = 3740.3333
$ws.Range("M116").Value = -298.3332999999998
$ws.Range("H122").Value = 3240.5186
$ws.Range("I122").Value = 3194.3809
$ws.Range("J122").Value = 3402
$ws.Range("K122").Value = 9583.1427
$ws.Range("L122").Value = 10206
$ws.Range("M122").Value = -7133.1427
$ws.Range("N122").Value = -15106
$ws.Range("H137").Value = 7083.625
$ws.Range("I137").Value = 5276.25
$ws.Range("K137").Value = 15828.75
$ws.Range("M137").Value = -13278.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2172952.5
$ws.Range("I32").Value = 1003354.06
$ws.Range("K32").Value = 1003354.06
$ws.Range("M32").Value = -1003067.06
$ws.Range("H61").Value = 4722
$ws.Range("J61").Value = 4937.5
$ws.Range("L61").Value = 4937.5
$ws.Range("N61").Value = -5361.5
$ws.Range("H135").Value = 108994.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 108994.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 108994.5
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -119134.5
$ws.Range("H136").Value = 4722
$ws.Range("J136").Value = 4937.5
$ws.Range("L136").Value = 14812.5
$ws.Range("N136").Value = -19912.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 62509544
$ws.Range("I20").Value = 125016296
$ws.Range("J20").Value = 2795
$ws.Range("K20").Value = 125016296
$ws.Range("L20").Value = 2795
$ws.Range("M20").Value = -125016049
$ws.Range("N20").Value = -3289
$ws.Range("H86").Value = 3375.3333
$ws.Range("I86").Value = 3480.3076
$ws.Range("K86").Value = 3480.3076
$ws.Range("M86").Value = -2357.3076
$ws.Range("H89").Value = 3375.3333
$ws.Range("I89").Value = 3480.3076
$ws.Range("K89").Value = 17401.538
$ws.Range("M89").Value = -11785.538
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("H134").Value = 1537.5
$ws.Range("I134").Value = 1103.2916
$ws.Range("J134").Value = 2579.6
$ws.Range("K134").Value = 3309.8748
$ws.Range("L134").Value = 7738.799999999999
$ws.Range("M134").Value = -774.8748000000001
$ws.Range("N134").Value = -12808.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 32746.295
$ws.Range("J4").Value = 38209.62
$ws.Range("L4").Value = 38209.62
$ws.Range("N4").Value = -38433.62
$ws.Range("H22").Value = 350
$ws.Range("J22").Value = 400
$ws.Range("L22").Value = 400
$ws.Range("N22").Value = -1100
$ws.Range("H99").Value = 2446
$ws.Range("J99").Value = 2880
$ws.Range("L99").Value = 2880
$ws.Range("N99").Value = -5876
$ws.Range("H126").Value = 2446
$ws.Range("J126").Value = 2880
$ws.Range("L126").Value = 8640
$ws.Range("N126").Value = -13580

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3928.1428
$ws.Range("I63").Value = 997
$ws.Range("J63").Value = 4416.6665
$ws.Range("K63").Value = 2991
$ws.Range("L63").Value = 13249.9995
$ws.Range("M63").Value = -2242
$ws.Range("N63").Value = -14747.9995
$ws.Range("H66").Value = 3928.1428
$ws.Range("I66").Value = 997
$ws.Range("J66").Value = 4416.6665
$ws.Range("K66").Value = 8973
$ws.Range("L66").Value = 39749.9985
$ws.Range("M66").Value = -5229
$ws.Range("N66").Value = -47237.9985
$ws.Range("H92").Value = 164
$ws.Range("J92").Value = 164
$ws.Range("L92").Value = 492
$ws.Range("N92").Value = -2988
$ws.Range("H131").Value = 11911456
$ws.Range("J131").Value = 1843.1538
$ws.Range("L131").Value = 5529.4614
$ws.Range("N131").Value = -15609.4614

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7141.7144
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 7141.7144
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 7141.7144
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -7681.7144
$ws.Range("H73").Value = 7141.7144
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 7141.7144
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 7141.7144
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -9013.714400000001
$ws.Range("H122").Value = 4138
$ws.Range("I122").Value = 3290.875
$ws.Range("K122").Value = 9872.625
$ws.Range("M122").Value = -7422.625
$ws.Range("H126").Value = 6325.1904
$ws.Range("I126").Value = 2113.25
$ws.Range("J126").Value = 11941.111
$ws.Range("K126").Value = 6339.75
$ws.Range("L126").Value = 35823.333
$ws.Range("M126").Value = -3869.75
$ws.Range("N126").Value = -40763.333
$ws.Range("H132").Value = 2380.6
$ws.Range("I132").Value = 2140.652
$ws.Range("K132").Value = 6421.956
$ws.Range("M132").Value = -3891.956
$ws.Range("H135").Value = 69999
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2525.375
$ws.Range("I7").Value = 2298.2
$ws.Range("J7").Value = 2904
$ws.Range("K7").Value = 2298.2
$ws.Range("L7").Value = 2904
$ws.Range("M7").Value = -2186.2
$ws.Range("N7").Value = -3128
$ws.Range("H22").Value = 119049500
$ws.Range("I22").Value = 10206223
$ws.Range("J22").Value = 500001000
$ws.Range("K22").Value = 10206223
$ws.Range("L22").Value = 500001000
$ws.Range("M22").Value = -10205928
$ws.Range("N22").Value = -500001590
$ws.Range("H27").Value = 119049500
$ws.Range("I27").Value = 10206223
$ws.Range("J27").Value = 500001000
$ws.Range("K27").Value = 10206223
$ws.Range("L27").Value = 500001000
$ws.Range("M27").Value = -10206116
$ws.Range("N27").Value = -500001214
$ws.Range("H55").Value = 543.55554
$ws.Range("I55").Value = 358.25
$ws.Range("J55").Value = 914.1667
$ws.Range("K55").Value = 358.25
$ws.Range("L55").Value = 914.1667
$ws.Range("M55").Value = -185.25
$ws.Range("N55").Value = -1260.1667
$ws.Range("H126").Value = 2525.375
$ws.Range("I126").Value = 2298.2
$ws.Range("J126").Value = 2904
$ws.Range("K126").Value = 6894.599999999999
$ws.Range("L126").Value = 8712
$ws.Range("M126").Value = -4424.599999999999
$ws.Range("N126").Value = -13652
$ws.Range("H137").Value = 57332.723
$ws.Range("J137").Value = 57764.117
$ws.Range("L137").Value = 57764.117
$ws.Range("N137").Value = -67964.117

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 27779998
$ws.Range("I122").Value = 2599.8
$ws.Range("J122").Value = 62501744
$ws.Range("K122").Value = 7799.400000000001
$ws.Range("L122").Value = 187505232
$ws.Range("M122").Value = -5349.400000000001
$ws.Range("N122").Value = -187510132

